$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Audi A4 B6 1.9 TDI 130 cv Gasóleo 2000 - 2004 AVF, AWX"
$ws.Range("B3").Value = "Peugeot 308 1 Serie 1.6 HDi 109 cv Gasóleo 2007 - 2014 9HY (DV6TED4), 9HZ (DV6TED4)"
$ws.Range("B4").Value = "Ford Fiesta 6 1.0 EcoBoost 100 cv Gasolina 2013 - 2024 SFJA, SFJB, SFJC, SFJD"
